$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (12 and 13) for the new error entries, pushing the
# existing 1000/1001/... rows down.
$ws.Rows("12:13").Insert()

# New row 13: Id 108 -> "手牌满惩罚" (written first so it claims the lower
# shared-string index, matching the authored sharedStrings.xml ordering)
$ws.Range("B13").Value = "手牌满惩罚"
$ws.Range("A13").Value = 108

# New row 12: Id 107 -> "卡牌耗尽惩罚"
$ws.Range("B12").Value = "卡牌耗尽惩罚"
$ws.Range("A12").Value = 107

# Resize the table (xl/tables/table1.xml) to cover the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B22"))

# Column A now gets an explicit width (matches the authored width="9").
$ws.Columns("A").ColumnWidth = 8.285714285714286

# Update the selected cell to match the post-edit selection.
$ws.Range("B12").Select()
